# Update for insert release-notes.md f80ed2bb9e1dd81abc71d13817b8a44a756cee80
#
# - Bump the Metadata table: Version, Status, Date, Contact
# - Swap the two "Mapping" columns (AK/AL) on the Elements sheet: header
#   text, per-row values, and column widths all move together.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value  = "0.4.0-snapshot-1"                 # Version
$meta.Range("B6").Value  = "draft"                            # Status
$meta.Range("B8").Value  = "2024-05-23T12:16:26+00:00"        # Date
$meta.Range("B10").Value = "ANS (https://esante.gouv.fr)"     # Contact

# --- Elements sheet: swap the AK (37) / AL (38) "Mapping" columns ---------
$elem = $wb.Worksheets.Item("Elements")

$lastRow = $elem.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $elem.Cells.Item($r, 37)
    $alCell = $elem.Cells.Item($r, 38)
    $akVal = $akCell.Value()
    $alVal = $alCell.Value()
    if ("$akVal" -ne "$alVal") {
        $akCell.Value = $alVal
        $alCell.Value = $akVal
    }
}

# Column widths follow the data: AK becomes the wide "business mapping"
# column, AL becomes the narrower "RIM Mapping" column.
$elem.Columns.Item(37).ColumnWidth = 83.57291666666667
$elem.Columns.Item(38).ColumnWidth = 24.147135416666668
